$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Populate the new columns E:G on the "Contrastive Pre-Trained" sheet
#    (this sheet gets renamed to "Class-Based Contrastive" afterwards).
#    Doing this BEFORE the rename / before touching T-Tests keeps the
#    shared-string insertion order (and therefore the new string
#    indices) lined up with the authored workbook.
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Contrastive Pre-Trained")

# Mirror the formatting of A1:C13 onto E1:G13 first (keeps the same
# striped-row style without disturbing any values), then fill in values.
$ws8.Range("A1:C13").Copy() | Out-Null
$ws8.Range("E1").PasteSpecial(-4122) | Out-Null

# Header row for the new block (mirrors A1:C1)
$ws8.Range("E1").Value = "Iteration"
$ws8.Range("F1").Value = "Dice Score"
$ws8.Range("G1").Value = "Hausdorff Distance"

# Iteration numbers
$ws8.Range("E2").Value = 1
$ws8.Range("E3").Value = 2
$ws8.Range("E4").Value = 3
$ws8.Range("E5").Value = 4
$ws8.Range("E6").Value = 5
$ws8.Range("E7").Value = 6
$ws8.Range("E8").Value = 7
$ws8.Range("E9").Value = 8
$ws8.Range("E10").Value = 9
$ws8.Range("E11").Value = 10

# Dice score / Hausdorff distance values (rows 6 and 10 are missing runs,
# left blank but already carry the correct formatting from the paste above)
$ws8.Range("F2").Value = 0.8166
$ws8.Range("G2").Value = 17.0245

$ws8.Range("F3").Value = 0.8534
$ws8.Range("G3").Value = 17.1648

$ws8.Range("F4").Value = 0.837
$ws8.Range("G4").Value = 15.6084

$ws8.Range("F5").Value = 0.8246
$ws8.Range("G5").Value = 15.5401

$ws8.Range("F7").Value = 0.8359
$ws8.Range("G7").Value = 17.4953

$ws8.Range("F8").Value = 0.8284
$ws8.Range("G8").Value = 15.0158

$ws8.Range("F9").Value = 0.8486
$ws8.Range("G9").Value = 17.2845

$ws8.Range("F11").Value = 0.8632
$ws8.Range("G11").Value = 18.0246

# Average / Standard Deviation rows
$ws8.Range("E12").Value = "Average"
$ws8.Range("F12").Formula = "=AVERAGE(F2:F11)"
$ws8.Range("G12").Formula = "=AVERAGE(G2:G11)"

$ws8.Range("E13").Value = "Standard Deviation"
$ws8.Range("F13").Formula = "=_xlfn.STDEV.S(F2:F11)"
$ws8.Range("G13").Formula = "=_xlfn.STDEV.S(G2:G11)"

# New note row describing the two blocks
$ws8.Range("A15").Value = "Full Training Set for fine-tuning"
$ws8.Range("E15").Value = "45% of training set for fine-tuning"

# Column widths for the new columns
$ws8.Columns.Item(5).ColumnWidth = 12.619791666666666
$ws8.Columns.Item(6).ColumnWidth = 11.256510416666666
$ws8.Columns.Item(7).ColumnWidth = 17.166666666666668

# ---------------------------------------------------------------------
# 2) Rename the sheet
# ---------------------------------------------------------------------
$ws8.Name = "Class-Based Contrastive"

# ---------------------------------------------------------------------
# 3) Update the T-Tests sheet: fix formulas referencing the renamed
#    sheet, and add the new "Reduced Supervised Contrastive vs Standard"
#    comparison block.
# ---------------------------------------------------------------------
$tt = $wb.Worksheets.Item("T-Tests")

$tt.Range("H2").Formula = "=_xlfn.T.TEST('Standard Training'!B2:B11, 'Class-Based Contrastive'!B2:B11, 2, 3)"
$tt.Range("H3").Formula = "=_xlfn.T.TEST('Standard Training'!C2:C11, 'Class-Based Contrastive'!C2:C11, 2, 3)"
$tt.Range("H6").Formula = "=_xlfn.T.TEST('Joint Training'!B2:B11, 'Class-Based Contrastive'!B2:B11, 2, 3)"
$tt.Range("H7").Formula = "=_xlfn.T.TEST('Joint Training'!C2:C11, 'Class-Based Contrastive'!C2:C11, 2, 3)"

$tt.Range("F13").Value = "Reduced Supervised Contrastive vs Standard"
$tt.Range("F14").Value = "Dice:"
$tt.Range("H14").Formula = "=_xlfn.T.TEST('Standard Training'!B2:B11, 'Class-Based Contrastive'!F2:F11, 2, 3)"
$tt.Range("F15").Value = "Hausdorff:"
$tt.Range("H15").Formula = "=_xlfn.T.TEST('Standard Training'!C2:C11, 'Class-Based Contrastive'!G2:G11, 2, 3)"

# ---------------------------------------------------------------------
# 4) Selection / active-tab bookkeeping to match the authored workbook
# ---------------------------------------------------------------------
$tt.Range("H16").Select() | Out-Null
$ws8.Range("F6").Select() | Out-Null
$ws8.Activate() | Out-Null
